# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" table (rows 16-35, columns B:G) is re-sorted by
# "Periodo Mora" (B..G = Tipo Doc, N Doc, Nombre, Periodo Mora, Valor Mora,
# Salario Basico) in ascending period order (1906 -> 2010), and the first
# periods (1909, 1911, 1912) for the new worker NORELVIS ESTHER ROSALES
# REDONDO are interleaved at their correct chronological position instead
# of being appended after JOHANNA's rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row data in B:G order: TipoDoc, NroDoc, Nombre, PeriodoMora, ValorMora, SalarioBasico
$johanna  = "CC", "1067899155", "JOHANNA MARIA BARRIOS COTTA"
$norelvis = "CC", "57461901",   "NORELVIS ESTHER ROSALES REDONDO"

$rows = @(
    @{ Row = 16; Worker = $johanna;  Periodo = "1906"; Mora = 23187; Salario = 828116 },
    @{ Row = 17; Worker = $johanna;  Periodo = "1907"; Mora = 33125; Salario = 828116 },
    @{ Row = 18; Worker = $johanna;  Periodo = "1908"; Mora = 33125; Salario = 828116 },
    @{ Row = 19; Worker = $johanna;  Periodo = "1909"; Mora = 33125; Salario = 828116 },
    @{ Row = 20; Worker = $norelvis; Periodo = "1909"; Mora = 48000; Salario = 1200000 },
    @{ Row = 21; Worker = $johanna;  Periodo = "1910"; Mora = 33125; Salario = 828116 },
    @{ Row = 22; Worker = $johanna;  Periodo = "1911"; Mora = 33125; Salario = 828116 },
    @{ Row = 23; Worker = $norelvis; Periodo = "1911"; Mora = 48000; Salario = 1200000 },
    @{ Row = 24; Worker = $johanna;  Periodo = "1912"; Mora = 33125; Salario = 828116 },
    @{ Row = 25; Worker = $norelvis; Periodo = "1912"; Mora = 48000; Salario = 1200000 },
    @{ Row = 26; Worker = $johanna;  Periodo = "2001"; Mora = 33125; Salario = 828116 },
    @{ Row = 27; Worker = $johanna;  Periodo = "2002"; Mora = 33125; Salario = 828116 },
    @{ Row = 28; Worker = $johanna;  Periodo = "2003"; Mora = 33125; Salario = 828116 },
    @{ Row = 29; Worker = $johanna;  Periodo = "2004"; Mora = 33125; Salario = 828116 },
    @{ Row = 30; Worker = $johanna;  Periodo = "2005"; Mora = 33125; Salario = 828116 },
    @{ Row = 31; Worker = $johanna;  Periodo = "2006"; Mora = 33125; Salario = 828116 },
    @{ Row = 32; Worker = $johanna;  Periodo = "2007"; Mora = 33125; Salario = 828116 },
    @{ Row = 33; Worker = $johanna;  Periodo = "2008"; Mora = 33125; Salario = 828116 },
    @{ Row = 34; Worker = $johanna;  Periodo = "2009"; Mora = 33125; Salario = 828116 },
    @{ Row = 35; Worker = $johanna;  Periodo = "2010"; Mora = 22083; Salario = 828116 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("B$rowNum").Value = $r.Worker[0]
    $ws.Range("C$rowNum").Value = $r.Worker[1]
    $ws.Range("D$rowNum").Value = $r.Worker[2]
    $ws.Range("E$rowNum").Value = $r.Periodo
    $ws.Range("F$rowNum").Value = $r.Mora
    $ws.Range("G$rowNum").Value = $r.Salario
}
